$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 3 into row 4 so the new row inherits the same
# per-cell styles (hyperlink style on C, wrap-text style on D) without
# registering brand-new style entries.
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)

# Populate the new row's values.
$ws.Range("A4").Value = "Arrays"
$ws.Range("B4").Value = "Remove Duplicates from Sorted Array"
$ws.Range("C4").Value = "https://leetcode.com/explore/interview/card/top-interview-questions-easy/92/array/727/"
$ws.Range("D4").Value = "Use two pointers :`nwhenever unique element is found store it at next index,for duplicates skip them."

# Attach the hyperlink for the question link cell.
$ws.Hyperlinks.Add($ws.Range("C4"), "https://leetcode.com/explore/interview/card/top-interview-questions-easy/92/array/727/")
$ws.Range("C4").Style = $ws.Range("C3").Style

# Row height for the new wrapped-text row.
$ws.Rows.Item(4).RowHeight = 43.5

# Column width adjustments (B/C) as a result of the new, longer content.
$ws.Columns.Item(2).ColumnWidth = 32.54296875
$ws.Columns.Item(3).ColumnWidth = 61.7265625

# Selection moves to A3:A4 with active cell A3.
$ws.Range("A3:A4").Select()
